# Table9_Master_Table_method.docx edit:
# Flip the section's page orientation from portrait to landscape.
# Word keeps the same physical page dimensions but swaps the
# width/height values and updates the orient attribute accordingly
# (16848 x 11952 twips portrait  ->  11952 x 16848 twips landscape).

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    $sec.PageSetup.Orientation = 1   # wdOrientLandscape
}
